$d = $word.ActiveDocument

# 1. Title: "Segmenting and Clustering Neighborhoods in New York City"
#    -> "Segmenting and Clustering Neighborhoods in Toronto City"
$d.Content.Find.Execute(
    "Segmenting and Clustering Neighborhoods in New York City",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Segmenting and Clustering Neighborhoods in Toronto City", 2)

# 2. "So, I would be interested in comparing the neighborhoods of the two cities
#    and determine how similar or dissimilar they are."
#    -> "... of the Toronto cities and determine ..."
$d.Content.Find.Execute(
    "interested in comparing the neighborhoods of the two cities",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "interested in comparing the neighborhoods of the Toronto cities", 2)

Write-Output "Done"
